$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.695.95'
$ws.Range("E2").Value = '  +2.51%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.904.18'
$ws.Range("E3").Value = '  +2.38%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.031'
$ws.Range("E4").Value = '  +2.69%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.10'
$ws.Range("E5").Value = '  +2.61%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.029'
$ws.Range("E6").Value = '  +2.59%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5187'
$ws.Range("E7").Value = '  +0.80%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3950'
$ws.Range("E8").Value = '  +3.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08341'
$ws.Range("E9").Value = '  +1.52%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.134'
$ws.Range("E10").Value = '  +2.11%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.49'
$ws.Range("E11").Value = '  +2.36%  '

$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.300'
$ws.Range("E12").Value = '  +1.77%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.913.89'
$ws.Range("E13").Value = '  +2.79%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.67'
$ws.Range("E14").Value = '  +0.53%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.314'
$ws.Range("E15").Value = '  +0.88%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.031'
$ws.Range("E16").Value = '  +2.78%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001113'
$ws.Range("E17").Value = '  +1.50%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.58'
$ws.Range("E18").Value = '  +1.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06794'
$ws.Range("E19").Value = '  +2.25%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.01'
$ws.Range("E20").Value = '  +1.82%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.028'
$ws.Range("E21").Value = '  +2.59%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.090'
$ws.Range("E22").Value = '  +1.39%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.720.75'
$ws.Range("E23").Value = '  +2.51%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.24'
$ws.Range("E24").Value = '  +1.59%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.281'
$ws.Range("E25").Value = '  +1.54%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.125.37'
$ws.Range("E26").Value = '  +2.48%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.39'
$ws.Range("E27").Value = '  +2.85%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.01'
$ws.Range("E28").Value = '  +2.76%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.462'
$ws.Range("E29").Value = '  -1.93%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.58'
$ws.Range("E30").Value = '  +2.34%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1060'
$ws.Range("E31").Value = '  -0.26%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.049'
$ws.Range("E32").Value = '  +1.83%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.006'
$ws.Range("E33").Value = '  +1.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.680'
$ws.Range("E34").Value = '  +2.21%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02469'
$ws.Range("E35").Value = '  +2.21%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.384'
$ws.Range("E36").Value = '  +0.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06616'
$ws.Range("E37").Value = '  +1.76%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2221'
$ws.Range("E38").Value = '  +2.14%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6557'
$ws.Range("E39").Value = '  +0.16%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.261'
$ws.Range("E40").Value = '  +3.67%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.195'
$ws.Range("E41").Value = '  +0.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.015'
$ws.Range("E42").Value = '  +0.48%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.16'
$ws.Range("E43").Value = '  -0.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6174'
$ws.Range("E44").Value = '  +0.31%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.18'
$ws.Range("E45").Value = '  +1.21%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.758'
$ws.Range("E46").Value = '  +2.44%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.301'
$ws.Range("E47").Value = '  +1.43%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.019'
$ws.Range("E48").Value = '  +0.78%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.240'
$ws.Range("E49").Value = '  +1.92%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '122.52'
$ws.Range("E50").Value = '  +1.63%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06959'
$ws.Range("E51").Value = '  +2.58%  '
